# Add new product rows to the "Productos" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

$ws.Range("A5").Value = "Americano"
$ws.Range("B5").Value = 4000.0
$ws.Range("C5").Value = 100.0

$ws.Range("A6").Value = "Latte"
$ws.Range("B6").Value = 5000.0
$ws.Range("C6").Value = 80.0

$ws.Range("A7").Value = "Capuccino"
$ws.Range("B7").Value = 6000.0
$ws.Range("C7").Value = 70.0

$ws.Range("A8").Value = "a"
$ws.Range("B8").Value = 1.0
$ws.Range("C8").Value = 2.0

$ws.Range("A9").Value = "a"
$ws.Range("B9").Value = 1.0
$ws.Range("C9").Value = 15.0
